$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Update header text for column G (shared string change)
$ws.Range("G1").Value = "Median_Liabilities"

# Row 2 (Error_Type id 7)
$ws.Range("B2").Value = 124608
$ws.Range("C2").Value = 0.269482442362897
$ws.Range("D2").Value = 0.248838880007541
$ws.Range("E2").Value = 0.249657076114419
$ws.Range("F2").Value = 0.230540219060229
$ws.Range("G2").Value = -0.0500788728053914

# Row 3 (Error_Type id 8)
$ws.Range("B3").Value = 44580
$ws.Range("C3").Value = -0.843690993736956
$ws.Range("D3").Value = -0.83887828567755
$ws.Range("E3").Value = -0.624553550835157
$ws.Range("F3").Value = -0.760497939050802
$ws.Range("G3").Value = 0.103074222886388

# Row 4 (Error_Type id 9)
$ws.Range("B4").Value = 175
$ws.Range("C4").Value = 0.0104436174992944
$ws.Range("D4").Value = -0.0400732388617184
$ws.Range("E4").Value = 0.0774335872267252
$ws.Range("F4").Value = 0.0290282858566791
$ws.Range("G4").Value = 0.0214540030466988

# Row 5 (Error_Type id 10)
$ws.Range("B5").Value = 1291
$ws.Range("C5").Value = -0.991959669159165
$ws.Range("D5").Value = -1.00821067376228
$ws.Range("E5").Value = -0.567911524566636
$ws.Range("F5").Value = -0.800272521643408
$ws.Range("G5").Value = 0.162951274448651
